# Automatische test-sync: 2025-08-03 23:31:50
# Appends the new "Testmail #3" row to the Logs sheet, extends the
# conditional-formatting ranges to cover it, and bumps the "Overig"
# tally on the Dashboard sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 52

$logs.Cells.Item($newRow, 1).Value = "Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #3: Hoi, hebben jullie al iets gehoord?"
$logs.Cells.Item($newRow, 4).Value = "Overig"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 23:31:12"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting sqref ranges from row 51 to row 52
# (one ModifyAppliesToRange per block is enough - it re-targets every
# cfRule sharing that sqref).
$logs.Range("D2:D51").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D52"))
$logs.Range("G2:G51").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G52"))
$logs.Range("H2:H51").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H52"))
$logs.Range("I2:I51").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I52"))
$logs.Range("J2:J51").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J52"))

# Bump the Dashboard "Overig" count (row 3) from 11 to 12.
$dashboard.Cells.Item(3, 2).Value = 12
